$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.400.22'
$ws.Cells.Item(2, 5).Value = '  +0.10%  '

$ws.Cells.Item(3, 4).Value = '3.510.99'
$ws.Cells.Item(3, 5).Value = '  +0.40%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '591.28'
$ws.Cells.Item(5, 5).Value = '  +0.89%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '134.59'
$ws.Cells.Item(6, 5).Value = '  -0.29%  '

$ws.Cells.Item(7, 5).Value = '  -0.04%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '7.64'
$ws.Cells.Item(9, 5).Value = '  +6.19%  '

$ws.Cells.Item(10, 5).Value = '  +0.60%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.389'
$ws.Cells.Item(11, 5).Value = '  +3.79%  '

$ws.Cells.Item(12, 4).Value = '4.108.22'
$ws.Cells.Item(12, 5).Value = '  +0.30%  '

$ws.Cells.Item(13, 5).Value = '  +1.10%  '

$ws.Cells.Item(14, 5).Value = '  +0.46%  '

$ws.Cells.Item(15, 4).Value = '3.509.73'
$ws.Cells.Item(15, 5).Value = '  +0.21%  '

$ws.Cells.Item(16, 2).Value = 'Avalanche'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '25.77'
$ws.Cells.Item(16, 5).Value = '  +2.97%  '

$ws.Cells.Item(17, 2).Value = 'WrappedBTC'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(17, 4).Value = '64.381.97'
$ws.Cells.Item(17, 5).Value = '  +0.01%  '

$ws.Cells.Item(18, 5).Value = '  -0.02%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '13.67'
$ws.Cells.Item(19, 5).Value = '  -0.61%  '

$ws.Cells.Item(20, 5).Value = '  +2.07%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '395.86'
$ws.Cells.Item(21, 5).Value = '  +2.92%  '

$ws.Cells.Item(22, 5).Value = '  +2.01%  '

$ws.Cells.Item(23, 4).Value = '3.650.12'
$ws.Cells.Item(23, 5).Value = '  +0.29%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '74.65'
$ws.Cells.Item(24, 5).Value = '  +0.80%  '

$ws.Cells.Item(25, 5).Value = '  +0.14%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '5.75'
$ws.Cells.Item(26, 5).Value = '  +0.47%  '

$ws.Cells.Item(27, 5).Value = '  +3.51%  '

$ws.Cells.Item(28, 5).Value = '  +0.07%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.41'
$ws.Cells.Item(29, 5).Value = '  -0.69%  '

$ws.Cells.Item(30, 5).Value = '  +1.20%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '8.24'
$ws.Cells.Item(31, 5).Value = '  +0.15%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.48'
$ws.Cells.Item(32, 5).Value = '  -4.08%  '

$ws.Cells.Item(33, 5).Value = '  +6.37%  '

$ws.Cells.Item(34, 4).Value = '3.538.85'
$ws.Cells.Item(34, 5).Value = '  +0.58%  '

$ws.Cells.Item(35, 5).Value = '  +0.03%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '23.38'
$ws.Cells.Item(36, 5).Value = '  -0.64%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '5.37'
$ws.Cells.Item(37, 5).Value = '  +1.70%  '

$ws.Cells.Item(38, 5).Value = '  +2.24%  '

$ws.Cells.Item(39, 5).Value = '  +1.11%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '167.36'
$ws.Cells.Item(40, 5).Value = '  +2.92%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0788'
$ws.Cells.Item(41, 5).Value = '  +1.01%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.812'
$ws.Cells.Item(42, 5).Value = '  +0.77%  '

$ws.Cells.Item(43, 5).Value = '  -0.02%  '

$ws.Cells.Item(44, 5).Value = '  +1.12%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '24.84'
$ws.Cells.Item(45, 5).Value = '  -3.78%  '

$ws.Cells.Item(46, 5).Value = '  +0.83%  '

$ws.Cells.Item(47, 5).Value = '  -2.68%  '

$ws.Cells.Item(48, 5).Value = '  +0.72%  '

$ws.Cells.Item(49, 4).Value = '2.377.27'
$ws.Cells.Item(49, 5).Value = '  -4.03%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.905'
$ws.Cells.Item(50, 5).Value = '  -0.35%  '

$ws.Cells.Item(51, 5).Value = '  +0.15%  '
